# Insert a new data row at row 55 (pushing the existing rows 55-108 down to 56-109)
# and populate it with the new record's values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(55).Insert()

$ws.Cells.Item(55, 1).Value = 10
$ws.Cells.Item(55, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(55, 3).Value = "La Araucanía"
$ws.Cells.Item(55, 4).Value = 44781
$ws.Cells.Item(55, 5).Value = 9
$ws.Cells.Item(55, 6).Value = 100112035
$ws.Cells.Item(55, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(55, 8).Value = "Sin especificar"
$ws.Cells.Item(55, 9).Value = "Primera"
$ws.Cells.Item(55, 10).Value = 130
$ws.Cells.Item(55, 11).Value = 25000
$ws.Cells.Item(55, 12).Value = 26000
$ws.Cells.Item(55, 13).Value = 25615
$ws.Cells.Item(55, 14).Value = "`$/malla 10 kilos"
$ws.Cells.Item(55, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(55, 16).Value = 2562
$ws.Cells.Item(55, 17).Value = 10
$ws.Cells.Item(55, 18).Value = "Hortaliza"
